# Combine RowFiller + ExcelMerger test fixture update:
#   - header cell B1 gets the "{name}" placeholder (new shared string)
#   - the sheet's active selection moves back to A1 (was A2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "{name}"
$ws.Range("A1").Select()
